$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Mensuração")
$ws.Columns("C:C").Insert()
